# Add the team's season record (Wins / Losses / Ties) as three new
# columns (AD, AE, AF) appended after the existing "Unnamed: 28" column (AC),
# for every player row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror the bold / centered / bordered header formatting used by the other
# column headers (row 1) by copying the format from the adjacent header cell.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1")) | Out-Null

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every data row (2-48): this team went 72-90-0.
$ws.Range("AD2:AD48").Value = 72
$ws.Range("AE2:AE48").Value = 90
$ws.Range("AF2:AF48").Value = 0
